# Sprint 2 Backlog update after scrum on day 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the in-progress items as "In progress" in the Status column (D)
$ws.Range("D3").Value = "In progress"
$ws.Range("D12").Value = "In progress"
$ws.Range("D17").Value = "In progress"
$ws.Range("D22").Value = "In progress"
$ws.Range("D27").Value = "In progress"

# Update the active selection to reflect where the user left off
$ws.Range("E25").Select()
